# "may revision w/o new cbo"
# Update monthly state UI figures for the May revision: a batch of cell-value
# corrections across existing rows (98-627), plus trailing revisions to
# rows 628 and a brand-new row 629 (2022-04-30) for columns A, P, Q.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(98, 2).Value = 2273261
$ws.Cells.Item(98, 6).Value = 13886960
$ws.Cells.Item(98, 7).Value = 11076395
$ws.Cells.Item(98, 8).Value = 904872
$ws.Cells.Item(98, 9).Value = 84.12000000000001
$ws.Cells.Item(98, 13).Value = 8643094
$ws.Cells.Item(104, 13).Value = 8353744
$ws.Cells.Item(108, 13).Value = 8263908
$ws.Cells.Item(132, 8).Value = 935342
$ws.Cells.Item(132, 9).Value = 100.69
$ws.Cells.Item(132, 13).Value = 14066889
$ws.Cells.Item(133, 13).Value = 14469341
$ws.Cells.Item(134, 13).Value = 14604977
$ws.Cells.Item(135, 13).Value = 14690413
$ws.Cells.Item(136, 13).Value = 14866873
$ws.Cells.Item(137, 13).Value = 14861391
$ws.Cells.Item(138, 13).Value = 14670782
$ws.Cells.Item(139, 13).Value = 14467007
$ws.Cells.Item(140, 13).Value = 14131787
$ws.Cells.Item(141, 13).Value = 13894084
$ws.Cells.Item(142, 13).Value = 13749990
$ws.Cells.Item(143, 13).Value = 13621859
$ws.Cells.Item(186, 4).Value = 590513
$ws.Cells.Item(186, 7).Value = 9819319
$ws.Cells.Item(186, 8).Value = 1209154
$ws.Cells.Item(186, 9).Value = 127.37
$ws.Cells.Item(186, 13).Value = 13889943
$ws.Cells.Item(186, 14).Value = 8211318
$ws.Cells.Item(187, 13).Value = 13966176
$ws.Cells.Item(187, 14).Value = 8231309
$ws.Cells.Item(188, 13).Value = 14177409
$ws.Cells.Item(188, 14).Value = 8311095
$ws.Cells.Item(189, 13).Value = 14260717
$ws.Cells.Item(189, 14).Value = 8315455
$ws.Cells.Item(190, 13).Value = 14423044
$ws.Cells.Item(190, 14).Value = 8371036
$ws.Cells.Item(191, 4).Value = 573542
$ws.Cells.Item(191, 5).Value = 207674
$ws.Cells.Item(191, 7).Value = 8670362
$ws.Cells.Item(191, 8).Value = 1073048
$ws.Cells.Item(191, 9).Value = 128.59
$ws.Cells.Item(191, 13).Value = 14528189
$ws.Cells.Item(191, 14).Value = 8368552
$ws.Cells.Item(191, 15).Value = 2551089
$ws.Cells.Item(192, 13).Value = 14532964
$ws.Cells.Item(192, 14).Value = 8290126
$ws.Cells.Item(192, 15).Value = 2547711
$ws.Cells.Item(193, 13).Value = 14761487
$ws.Cells.Item(193, 14).Value = 8365218
$ws.Cells.Item(193, 15).Value = 2572359
$ws.Cells.Item(194, 6).Value = 14729048
$ws.Cells.Item(194, 13).Value = 14844168
$ws.Cells.Item(194, 14).Value = 8388105
$ws.Cells.Item(194, 15).Value = 2570159
$ws.Cells.Item(195, 13).Value = 14889099
$ws.Cells.Item(195, 14).Value = 8325577
$ws.Cells.Item(195, 15).Value = 2571048
$ws.Cells.Item(196, 6).Value = 13097775
$ws.Cells.Item(196, 13).Value = 14967780
$ws.Cells.Item(196, 14).Value = 8320073
$ws.Cells.Item(196, 15).Value = 2563122
$ws.Cells.Item(197, 13).Value = 15089925
$ws.Cells.Item(197, 14).Value = 8341537
$ws.Cells.Item(197, 15).Value = 2566577
$ws.Cells.Item(198, 13).Value = 15141496
$ws.Cells.Item(198, 14).Value = 8330302
$ws.Cells.Item(198, 15).Value = 2567378
$ws.Cells.Item(199, 5).Value = 214630
$ws.Cells.Item(199, 6).Value = 10187328
$ws.Cells.Item(199, 13).Value = 15309040
$ws.Cells.Item(199, 14).Value = 8356550
$ws.Cells.Item(199, 15).Value = 2587932
$ws.Cells.Item(200, 13).Value = 15440650
$ws.Cells.Item(200, 14).Value = 8350569
$ws.Cells.Item(200, 15).Value = 2599986
$ws.Cells.Item(201, 13).Value = 15510424
$ws.Cells.Item(201, 14).Value = 8368312
$ws.Cells.Item(201, 15).Value = 2605564
$ws.Cells.Item(202, 6).Value = 10071201
$ws.Cells.Item(202, 13).Value = 15673734
$ws.Cells.Item(202, 14).Value = 8381869
$ws.Cells.Item(202, 15).Value = 2635477
$ws.Cells.Item(203, 6).Value = 9851485
$ws.Cells.Item(203, 15).Value = 2650976
$ws.Cells.Item(204, 15).Value = 2658798
$ws.Cells.Item(205, 15).Value = 2687347
$ws.Cells.Item(206, 15).Value = 2682401
$ws.Cells.Item(207, 15).Value = 2694752
$ws.Cells.Item(208, 15).Value = 2719335
$ws.Cells.Item(209, 15).Value = 2707777
$ws.Cells.Item(210, 15).Value = 2680342
$ws.Cells.Item(232, 2).Value = 1361133
$ws.Cells.Item(232, 6).Value = 11589001
$ws.Cells.Item(495, 7).Value = 15692944
$ws.Cells.Item(496, 7).Value = 16788728
$ws.Cells.Item(497, 7).Value = 13414873
$ws.Cells.Item(530, 2).Value = 1861804
$ws.Cells.Item(530, 6).Value = 14351206
$ws.Cells.Item(587, 2).Value = 899689
$ws.Cells.Item(587, 7).Value = 5490202
$ws.Cells.Item(588, 2).Value = 1052002
$ws.Cells.Item(588, 7).Value = 5165204
$ws.Cells.Item(589, 2).Value = 1229360
$ws.Cells.Item(589, 7).Value = 6603982
$ws.Cells.Item(589, 9).Value = 363.19
$ws.Cells.Item(626, 2).Value = 1413899
$ws.Cells.Item(626, 6).Value = 9145054
$ws.Cells.Item(626, 7).Value = 7008256
$ws.Cells.Item(626, 9).Value = 390.26
$ws.Cells.Item(626, 18).Value = 323638
$ws.Cells.Item(627, 2).Value = 905379
$ws.Cells.Item(627, 4).Value = 391606
$ws.Cells.Item(627, 5).Value = 139258
$ws.Cells.Item(627, 6).Value = 7687791
$ws.Cells.Item(627, 7).Value = 6243367
$ws.Cells.Item(627, 8).Value = 2411700
$ws.Cells.Item(627, 9).Value = 399.95
$ws.Cells.Item(627, 12).Value = 355.83
$ws.Cells.Item(627, 13).Value = 38621166
$ws.Cells.Item(627, 14).Value = 6871547
$ws.Cells.Item(627, 15).Value = 3605520
$ws.Cells.Item(627, 16).Value = 5212
$ws.Cells.Item(627, 17).Value = 14073
$ws.Cells.Item(627, 18).Value = 323990
$ws.Cells.Item(628, 12).Value = 359.74
$ws.Cells.Item(628, 13).Value = 35811117
$ws.Cells.Item(628, 14).Value = 6391307
$ws.Cells.Item(628, 15).Value = 3165310
$ws.Cells.Item(628, 16).Value = 5200
$ws.Cells.Item(628, 17).Value = 14091
$ws.Cells.Item(628, 18).Value = 322457
$ws.Cells.Item(629, 1).Value = 44681
$ws.Cells.Item(629, 16).Value = 5207
$ws.Cells.Item(629, 17).Value = 14112
